$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8031004071235657
$ws.Range("B1").Value = 1.509137511253357
$ws.Range("C1").Value = 5.876585483551025
$ws.Range("D1").Value = 3.124348640441895
$ws.Range("E1").Value = 1.466284513473511
